# Applies the "Automatic update of files" diff to the PARTILLE overview sheet.
#
# Net effect of the diff:
#   - Column C ("Förändrad") bumps from 46064 to 46065 for every data row (2-9).
#   - Rows 4-8 get reshuffled: the (A/B/F/G) tuple that used to sit on one row
#     now sits on another row (D/E and H..R are identical across these rows,
#     so they are unaffected). We apply the resulting final values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: "Förändrad" date serial, bumped by one day for every row ---
$ws.Range("C2").Value = 46065
$ws.Range("C3").Value = 46065
$ws.Range("C4").Value = 46065
$ws.Range("C5").Value = 46065
$ws.Range("C6").Value = 46065
$ws.Range("C7").Value = 46065
$ws.Range("C8").Value = 46065
$ws.Range("C9").Value = 46065

# --- Rows 4-8: A (Beteckning), B (Datum), F (Markägare), G (Area) reshuffled ---

# Row 4 -> becomes old row 8's data; no "Markägare" (F) entry
$ws.Range("A4").Value = "A 35734-2023"
$ws.Range("B4").Value = 45147.89258101852
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = 5.9

# Row 5 -> becomes old row 7's data; no "Markägare" (F) entry
$ws.Range("A5").Value = "A 26074-2025"
$ws.Range("B5").Value = 45805.32366898148
$ws.Range("G5").Value = 1.3

# Row 6 -> becomes old row 4's data; keeps "Kommuner" in F
$ws.Range("A6").Value = "A 25251-2025"
$ws.Range("B6").Value = 45800.50082175926
$ws.Range("F6").Value = "Kommuner"
$ws.Range("G6").Value = 0.7

# Row 7 -> becomes old row 5's data; no "Markägare" (F) entry
$ws.Range("A7").Value = "A 6983-2023"
$ws.Range("B7").Value = 44967.68585648148
$ws.Range("G7").Value = 5.4

# Row 8 -> becomes old row 6's data; gains "Kommuner" in F
$ws.Range("A8").Value = "A 25254-2025"
$ws.Range("B8").Value = 45800.50479166667
$ws.Range("F8").Value = "Kommuner"
$ws.Range("G8").Value = 0.2
